$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto price/volume refresh.
# Price cells that look like plain decimals (e.g. "1.000", "0.5120")
# are forced to Text format before assignment so Excel keeps the exact
# digits/trailing zeros instead of silently converting them to numbers.

$ws.Range("D2").Value = "26.826.46"
$ws.Range("E2").Value = "  +4.34%  "

$ws.Range("D3").Value = "1.875.25"
$ws.Range("E3").Value = "  +3.22%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.14"
$ws.Range("E5").Value = "  -0.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5268"
$ws.Range("E7").Value = "  +3.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3407"
$ws.Range("E8").Value = "  -3.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06946"
$ws.Range("E9").Value = "  +4.31%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.06"
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8041"
$ws.Range("E11").Value = "  -2.54%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07722"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").Value = "1.836.77"
$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.183"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.27"
$ws.Range("E15").Value = "  +3.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  +3.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008040"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "26.856.02"
$ws.Range("E20").Value = "  +4.24%  "

$ws.Range("D21").Value = "2.079.88"
$ws.Range("E21").Value = "  +0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.746"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("E23").Value = "  +0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.175"
$ws.Range("E24").Value = "  +1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.392"
$ws.Range("E25").Value = "  +8.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.39"
$ws.Range("E26").Value = "  +2.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.33"

$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.52"
$ws.Range("E29").Value = "  +3.59%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.356"
$ws.Range("E30").Value = "  +0.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.302"
$ws.Range("E31").Value = "  +1.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08906"
$ws.Range("E32").Value = "  +1.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04912"
$ws.Range("E33").Value = "  +0.67%  "

$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("E35").Value = "  -0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.871"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.275"
$ws.Range("E37").Value = "  +4.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.344"
$ws.Range("E38").Value = "  -0.91%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01847"
$ws.Range("E39").Value = "  -0.17%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5126"
$ws.Range("E40").Value = "  -0.74%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9567"
$ws.Range("E41").Value = "  -0.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "116.45"
$ws.Range("E42").Value = "  +5.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.157"
$ws.Range("E43").Value = "  -1.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.125"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4461"
$ws.Range("E46").Value = "  -2.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1336"
$ws.Range("E47").Value = "  -1.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.353"
$ws.Range("E48").Value = "  +1.50%  "

$ws.Range("E49").Value = "  -0.86%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05940"
$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.493"
$ws.Range("E51").Value = "  -0.37%  "

